$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H137").Value = 5717027
$ws.Range("I137").Value = 1966.7142
$ws.Range("J137").Value = 14289617
$ws.Range("K137").Value = 5900.142599999999
$ws.Range("L137").Value = 42868851
$ws.Range("M137").Value = -3350.142599999999
$ws.Range("N137").Value = -42873951

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H138").Value = 5816790.5
$ws.Range("I138").Value = 1693.3334
$ws.Range("J138").Value = 13162176
$ws.Range("K138").Value = 5080.0002
$ws.Range("L138").Value = 39486528
$ws.Range("M138").Value = 59.9997999999996
$ws.Range("N138").Value = -39496808

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H34").Value = 21776.8
$ws.Range("J34").Value = 21776.8
$ws.Range("L34").Value = 21776.8
$ws.Range("N34").Value = -22318.8

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H110").Value = 1620.4
$ws.Range("I110").Value = 1481.8334
$ws.Range("J110").Value = 1828.25
$ws.Range("K110").Value = 1481.8334
$ws.Range("L110").Value = 1828.25
$ws.Range("M110").Value = 563.1666
$ws.Range("N110").Value = -5918.25

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H132").Value = 20838470
$ws.Range("I132").Value = 50006404
$ws.Range("K132").Value = 150019212
$ws.Range("M132").Value = -150016682

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H107").Value = 2265.3809
$ws.Range("I107").Value = 2625.6428
$ws.Range("J107").Value = 1544.8572
$ws.Range("K107").Value = 2625.6428
$ws.Range("L107").Value = 1544.8572
$ws.Range("M107").Value = -705.6428000000001
$ws.Range("N107").Value = -5384.8572

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H4").Value = 76175710
$ws.Range("I4").Value = 167350000
$ws.Range("K4").Value = 167350000
$ws.Range("M4").Value = -167349888

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 7580270.5
$ws.Range("I31").Value = 7008.08
$ws.Range("K31").Value = 7008.08
$ws.Range("M31").Value = -6713.08

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H34").Value = 7580270.5
$ws.Range("I34").Value = 7008.08
$ws.Range("K34").Value = 7008.08
$ws.Range("M34").Value = -6806.08

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H4").Value = 5000170
$ws.Range("I4").Value = 5000170
$ws.Range("K4").Value = 15000510
$ws.Range("M4").Value = -15000398

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H107").Value = 880.0599999999999
$ws.Range("I107").Value = 321.16
$ws.Range("J107").Value = 1438.96
$ws.Range("K107").Value = 963.48
$ws.Range("L107").Value = 4316.88
$ws.Range("M107").Value = 956.52
$ws.Range("N107").Value = -8156.88

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H121").Value = 1083.25
$ws.Range("I121").Value = 400
$ws.Range("K121").Value = 1200
$ws.Range("M121").Value = 110

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H122").Value = 1393.0834
$ws.Range("I122").Value = 5142
$ws.Range("J122").Value = 643.3
$ws.Range("K122").Value = 46278
$ws.Range("L122").Value = 5789.7
$ws.Range("M122").Value = -43828
$ws.Range("N122").Value = -10689.7

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H141").Value = 2132.8572
$ws.Range("I141").Value = 2132.8572
$ws.Range("K141").Value = 6398.571599999999
$ws.Range("M141").Value = -1218.571599999999

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H5").Value = 500
$ws.Range("I5").Value = 500
$ws.Range("K5").Value = 500
$ws.Range("M5").Value = -388

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H33").Value = 9973.333000000001
$ws.Range("I33").Value = 9954
$ws.Range("J33").Value = 9997.5
$ws.Range("K33").Value = 9954
$ws.Range("L33").Value = 9997.5
$ws.Range("M33").Value = -9702
$ws.Range("N33").Value = -10501.5

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H62").Value = 0
$ws.Range("J62").Value = 0
$ws.Range("L62").Value = 0
$ws.Range("N62").Value = ""

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H65").Value = 0
$ws.Range("J65").Value = 0
$ws.Range("L65").Value = 0
$ws.Range("N65").Value = ""

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H97").Value = 338.6154
$ws.Range("I97").Value = 428
$ws.Range("J97").Value = 137.5
$ws.Range("K97").Value = 428
$ws.Range("L97").Value = 137.5
$ws.Range("M97").Value = 68
$ws.Range("N97").Value = -1129.5

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H132").Value = 4066.848
$ws.Range("I132").Value = 3041.7188
$ws.Range("J132").Value = 6410
$ws.Range("K132").Value = 9125.1564
$ws.Range("L132").Value = 19230
$ws.Range("M132").Value = -6595.1564
$ws.Range("N132").Value = -24290

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H2").Value = 49592.332
$ws.Range("J2").Value = 49592.332
$ws.Range("L2").Value = 49592.332
$ws.Range("N2").Value = -49816.332

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 905.55554
$ws.Range("I22").Value = 837.625
$ws.Range("J22").Value = 959.9
$ws.Range("K22").Value = 837.625
$ws.Range("L22").Value = 959.9
$ws.Range("M22").Value = -542.625
$ws.Range("N22").Value = -1549.9

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H27").Value = 905.55554
$ws.Range("I27").Value = 837.625
$ws.Range("J27").Value = 959.9
$ws.Range("K27").Value = 837.625
$ws.Range("L27").Value = 959.9
$ws.Range("M27").Value = -730.625
$ws.Range("N27").Value = -1173.9

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H41").Value = 189407.4
$ws.Range("I41").Value = 230000
$ws.Range("J41").Value = 27037
$ws.Range("K41").Value = 230000
$ws.Range("L41").Value = 27037
$ws.Range("M41").Value = -229562
$ws.Range("N41").Value = -27913

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H43").Value = 27014
$ws.Range("J43").Value = 27014
$ws.Range("L43").Value = 27014
$ws.Range("N43").Value = -27400

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H63").Value = 20085
$ws.Range("J63").Value = 20085
$ws.Range("L63").Value = 20085
$ws.Range("N63").Value = -21583

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H66").Value = 20085
$ws.Range("J66").Value = 20085
$ws.Range("L66").Value = 60255
$ws.Range("N66").Value = -67743

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H122").Value = 6595.8096
$ws.Range("I122").Value = 8805.5
$ws.Range("K122").Value = 26416.5
$ws.Range("M122").Value = -23966.5

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H33").Value = 16192.1
$ws.Range("J33").Value = 16192.1
$ws.Range("L33").Value = 16192.1
$ws.Range("N33").Value = -16692.1

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H36").Value = 16192.1
$ws.Range("J36").Value = 16192.1
$ws.Range("L36").Value = 16192.1
$ws.Range("N36").Value = -16692.1

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H37").Value = 0
$ws.Range("J37").Value = 0
$ws.Range("L37").Value = 0
$ws.Range("N37").Value = ""

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H42").Value = 24186
$ws.Range("J42").Value = 26566.666
$ws.Range("L42").Value = 26566.666
$ws.Range("N42").Value = -27322.666

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H43").Value = 14009
$ws.Range("I43").Value = 14009
$ws.Range("J43").Value = 0
$ws.Range("K43").Value = 14009
$ws.Range("L43").Value = 0
$ws.Range("M43").Value = -13860
$ws.Range("N43").Value = ""

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H100").Value = 1350
$ws.Range("I100").Value = 800
$ws.Range("J100").Value = 1460
$ws.Range("K100").Value = 1600
$ws.Range("L100").Value = 2920
$ws.Range("M100").Value = -1059
$ws.Range("N100").Value = -4002

Write-Host "All updates applied"